$wb = $excel.ActiveWorkbook

# ===== Sheet: LP1912 =====
$ws = $wb.Worksheets.Item("LP1912")

# Row 2: update A
$ws.Cells.Item(2,1).Value = "Última actualización: 09:28:08"
# Row 3: update A
$ws.Cells.Item(3,1).Value = "Total filas: 148"
# Row 26: update A,C,D
$ws.Cells.Item(26,1).Value = "06:25:28"
$ws.Cells.Item(26,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(26,4).Value = 5
# Row 27: update A,C,D
$ws.Cells.Item(27,1).Value = "05:55:02"
$ws.Cells.Item(27,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(27,4).Value = 35
# Row 32: update C
$ws.Cells.Item(32,3).Value = "215C_EL PATO"
# Row 33: update C
$ws.Cells.Item(33,3).Value = "14_ABASTO"
# Row 49: update C
$ws.Cells.Item(49,3).Value = "11_ETCHEVERRY"
# Row 51: update C
$ws.Cells.Item(51,3).Value = "84_COLONIA URQUIZA-ESC 49"
# Row 52: update C
$ws.Cells.Item(52,3).Value = "84_COLONIA URQUIZA-ESC 49"
# Row 53: update A,C,D
$ws.Cells.Item(53,1).Value = "05:55:02"
$ws.Cells.Item(53,3).Value = "16_SANTA ANA"
$ws.Cells.Item(53,4).Value = 97
# Row 54: update A,C,D
$ws.Cells.Item(54,1).Value = "06:54:06"
$ws.Cells.Item(54,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(54,4).Value = 38
# Row 78: update C
$ws.Cells.Item(78,3).Value = "215B_EL PATO"
# Row 79: update C
$ws.Cells.Item(79,3).Value = "16_P MOR-SANTA ANA"
# Row 80: update C
$ws.Cells.Item(80,3).Value = "215B_EL PATO"
# Row 81: update C
$ws.Cells.Item(81,3).Value = "16_P MOR-SANTA ANA"
# Row 93: update C
$ws.Cells.Item(93,3).Value = "16_SANTA ANA"
# Row 94: update C
$ws.Cells.Item(94,3).Value = "14_ABASTO"
# Row 97: update A,C,D
$ws.Cells.Item(97,1).Value = "07:17:59"
$ws.Cells.Item(97,3).Value = "17_ROMERO"
$ws.Cells.Item(97,4).Value = 96
# Row 98: update A,C,D
$ws.Cells.Item(98,1).Value = "08:47:26"
$ws.Cells.Item(98,3).Value = "10_OLMOS"
$ws.Cells.Item(98,4).Value = 6
# Row 114: update A,C,D
$ws.Cells.Item(114,1).Value = "08:47:26"
$ws.Cells.Item(114,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(114,4).Value = 36
# Row 116: update A,C,D
$ws.Cells.Item(116,1).Value = "08:55:01"
$ws.Cells.Item(116,3).Value = "16_SANTA ANA"
$ws.Cells.Item(116,4).Value = 28
# Row 118: update A,B,C,D
$ws.Cells.Item(118,1).Value = "09:28:08"
$ws.Cells.Item(118,2).Value = "09:29"
$ws.Cells.Item(118,3).Value = "17_ROMERO"
$ws.Cells.Item(118,4).Value = 1
# Row 119: update A,B,C,D
$ws.Cells.Item(119,1).Value = "07:48:05"
$ws.Cells.Item(119,2).Value = "09:32"
$ws.Cells.Item(119,3).Value = "15_ABASTO"
$ws.Cells.Item(119,4).Value = 104
# Row 120: update A,B,C,D
$ws.Cells.Item(120,1).Value = "09:28:08"
$ws.Cells.Item(120,2).Value = "09:33"
$ws.Cells.Item(120,3).Value = "10_OLMOS"
$ws.Cells.Item(120,4).Value = 5
# Row 121: update A,C,D
$ws.Cells.Item(121,1).Value = "07:48:05"
$ws.Cells.Item(121,3).Value = "16_SANTA ANA"
$ws.Cells.Item(121,4).Value = 106
# Row 122: update A,B,D
$ws.Cells.Item(122,1).Value = "08:47:26"
$ws.Cells.Item(122,2).Value = "09:34"
$ws.Cells.Item(122,4).Value = 47
# Row 123: update A,C,D
$ws.Cells.Item(123,1).Value = "08:55:01"
$ws.Cells.Item(123,3).Value = "16_SANTA ANA"
$ws.Cells.Item(123,4).Value = 40
# Row 124: update A,C,D
$ws.Cells.Item(124,1).Value = "08:47:26"
$ws.Cells.Item(124,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(124,4).Value = 48
# Row 125: update A,B,C,D
$ws.Cells.Item(125,1).Value = "09:28:08"
$ws.Cells.Item(125,2).Value = "09:35"
$ws.Cells.Item(125,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(125,4).Value = 7
# Row 126: update A,B,C,D
$ws.Cells.Item(126,1).Value = "09:28:08"
$ws.Cells.Item(126,2).Value = "09:42"
$ws.Cells.Item(126,3).Value = "215C_EL PATO"
$ws.Cells.Item(126,4).Value = 14
# Row 127: update A,C,D
$ws.Cells.Item(127,1).Value = "08:55:01"
$ws.Cells.Item(127,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(127,4).Value = 48
# Row 128: update A,B,D
$ws.Cells.Item(128,1).Value = "08:47:26"
$ws.Cells.Item(128,2).Value = "09:43"
$ws.Cells.Item(128,4).Value = 56
# Row 129: update A,B,C,D
$ws.Cells.Item(129,1).Value = "09:28:08"
$ws.Cells.Item(129,2).Value = "09:44"
$ws.Cells.Item(129,3).Value = "14_ABASTO"
$ws.Cells.Item(129,4).Value = 16
# Row 130: update A,B,D
$ws.Cells.Item(130,1).Value = "09:28:08"
$ws.Cells.Item(130,2).Value = "09:51"
$ws.Cells.Item(130,4).Value = 23
# Row 131: update A,B,C,D
$ws.Cells.Item(131,1).Value = "09:28:08"
$ws.Cells.Item(131,2).Value = "09:52"
$ws.Cells.Item(131,3).Value = "15_ABASTO"
$ws.Cells.Item(131,4).Value = 24
# Row 132: update A,B,C,D
$ws.Cells.Item(132,1).Value = "09:28:08"
$ws.Cells.Item(132,2).Value = "09:53"
$ws.Cells.Item(132,3).Value = "10_OLMOS"
$ws.Cells.Item(132,4).Value = 25
# Row 133: update A,B,C,D
$ws.Cells.Item(133,1).Value = "08:31:01"
$ws.Cells.Item(133,2).Value = "09:56"
$ws.Cells.Item(133,3).Value = "10_OLMOS"
$ws.Cells.Item(133,4).Value = 85
# Row 134: update A,B,C,D
$ws.Cells.Item(134,1).Value = "09:28:08"
$ws.Cells.Item(134,2).Value = "10:04"
$ws.Cells.Item(134,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(134,4).Value = 36
# Row 135: update A,B,C,D
$ws.Cells.Item(135,1).Value = "09:28:08"
$ws.Cells.Item(135,2).Value = "10:05"
$ws.Cells.Item(135,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(135,4).Value = 37
# Row 136: update A,B,C,D
$ws.Cells.Item(136,1).Value = "08:47:26"
$ws.Cells.Item(136,2).Value = "10:10"
$ws.Cells.Item(136,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(136,4).Value = 83
# Row 137: update A,B,C,D
$ws.Cells.Item(137,1).Value = "09:28:08"
$ws.Cells.Item(137,2).Value = "10:11"
$ws.Cells.Item(137,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(137,4).Value = 43
# Row 138: update A,B,C,D
$ws.Cells.Item(138,1).Value = "09:28:08"
$ws.Cells.Item(138,2).Value = "10:12"
$ws.Cells.Item(138,3).Value = "15_ABASTO"
$ws.Cells.Item(138,4).Value = 44
# Row 139: update A,B,C,D
$ws.Cells.Item(139,1).Value = "09:28:08"
$ws.Cells.Item(139,2).Value = "10:16"
$ws.Cells.Item(139,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(139,4).Value = 48
# Row 140: update A,B,C,D
$ws.Cells.Item(140,1).Value = "09:28:08"
$ws.Cells.Item(140,2).Value = "10:21"
$ws.Cells.Item(140,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(140,4).Value = 53
# New row 141
$ws.Cells.Item(141,1).Value = "09:28:08"
$ws.Cells.Item(141,2).Value = "10:23"
$ws.Cells.Item(141,3).Value = "16_SANTA ANA"
$ws.Cells.Item(141,4).Value = 55
$ws.Cells.Item(141,5).Value = "LP1912"
# New row 142
$ws.Cells.Item(142,1).Value = "09:28:08"
$ws.Cells.Item(142,2).Value = "10:24"
$ws.Cells.Item(142,3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(142,4).Value = 56
$ws.Cells.Item(142,5).Value = "LP1912"
# New row 143
$ws.Cells.Item(143,1).Value = "08:55:01"
$ws.Cells.Item(143,2).Value = "10:26"
$ws.Cells.Item(143,3).Value = "10_OLMOS"
$ws.Cells.Item(143,4).Value = 91
$ws.Cells.Item(143,5).Value = "LP1912"
# New row 144
$ws.Cells.Item(144,1).Value = "08:47:26"
$ws.Cells.Item(144,2).Value = "10:26"
$ws.Cells.Item(144,3).Value = "215A_EL PATO"
$ws.Cells.Item(144,4).Value = 99
$ws.Cells.Item(144,5).Value = "LP1912"
# New row 145
$ws.Cells.Item(145,1).Value = "09:28:08"
$ws.Cells.Item(145,2).Value = "10:27"
$ws.Cells.Item(145,3).Value = "215A_EL PATO"
$ws.Cells.Item(145,4).Value = 59
$ws.Cells.Item(145,5).Value = "LP1912"
# New row 146
$ws.Cells.Item(146,1).Value = "09:28:08"
$ws.Cells.Item(146,2).Value = "10:35"
$ws.Cells.Item(146,3).Value = "16_SANTA ANA"
$ws.Cells.Item(146,4).Value = 67
$ws.Cells.Item(146,5).Value = "LP1912"
# New row 147
$ws.Cells.Item(147,1).Value = "09:28:08"
$ws.Cells.Item(147,2).Value = "10:42"
$ws.Cells.Item(147,3).Value = "17_ROMERO"
$ws.Cells.Item(147,4).Value = 74
$ws.Cells.Item(147,5).Value = "LP1912"
# New row 148
$ws.Cells.Item(148,1).Value = "08:47:26"
$ws.Cells.Item(148,2).Value = "10:43"
$ws.Cells.Item(148,3).Value = "14_ABASTO"
$ws.Cells.Item(148,4).Value = 116
$ws.Cells.Item(148,5).Value = "LP1912"
# New row 149
$ws.Cells.Item(149,1).Value = "09:28:08"
$ws.Cells.Item(149,2).Value = "10:44"
$ws.Cells.Item(149,3).Value = "14_ABASTO"
$ws.Cells.Item(149,4).Value = 76
$ws.Cells.Item(149,5).Value = "LP1912"
# New row 150
$ws.Cells.Item(150,1).Value = "09:28:08"
$ws.Cells.Item(150,2).Value = "11:01"
$ws.Cells.Item(150,3).Value = "27_EL RETIRO"
$ws.Cells.Item(150,4).Value = 93
$ws.Cells.Item(150,5).Value = "LP1912"
# New row 151
$ws.Cells.Item(151,1).Value = "09:28:08"
$ws.Cells.Item(151,2).Value = "11:02"
$ws.Cells.Item(151,3).Value = "215C_EL PATO"
$ws.Cells.Item(151,4).Value = 94
$ws.Cells.Item(151,5).Value = "LP1912"
# New row 152
$ws.Cells.Item(152,1).Value = "09:28:08"
$ws.Cells.Item(152,2).Value = "11:07"
$ws.Cells.Item(152,3).Value = "16_P MOR-167 Y 521"
$ws.Cells.Item(152,4).Value = 99
$ws.Cells.Item(152,5).Value = "LP1912"
# New row 153
$ws.Cells.Item(153,1).Value = "09:28:08"
$ws.Cells.Item(153,2).Value = "11:20"
$ws.Cells.Item(153,3).Value = "86_EST CHICA-ESC AGRARIA"
$ws.Cells.Item(153,4).Value = 112
$ws.Cells.Item(153,5).Value = "LP1912"

# ===== Sheet: LP1912-215 =====
$ws = $wb.Worksheets.Item("LP1912-215")

# Row 2: update A
$ws.Cells.Item(2,1).Value = "Última actualización: 09:28:08"
# Row 3: update A
$ws.Cells.Item(3,1).Value = "Total filas: 21"
# Row 23: update A,D
$ws.Cells.Item(23,1).Value = "09:28:08"
$ws.Cells.Item(23,4).Value = 14
# Row 25: update A,D
$ws.Cells.Item(25,1).Value = "09:28:08"
$ws.Cells.Item(25,4).Value = 59
# New row 26
$ws.Cells.Item(26,1).Value = "09:28:08"
$ws.Cells.Item(26,2).Value = "11:02"
$ws.Cells.Item(26,3).Value = "215C_EL PATO"
$ws.Cells.Item(26,4).Value = 94
$ws.Cells.Item(26,5).Value = "LP1912"

# ===== Sheet: 6203-6173 =====
$ws = $wb.Worksheets.Item("6203-6173")

# Row 2: update A
$ws.Cells.Item(2,1).Value = "Última actualización: 09:28:08"
# Row 3: update A
$ws.Cells.Item(3,1).Value = "Total filas: 23"
# Row 26: update A,D
$ws.Cells.Item(26,1).Value = "09:28:08"
$ws.Cells.Item(26,4).Value = 35
# Row 27: update A,D
$ws.Cells.Item(27,1).Value = "09:28:08"
$ws.Cells.Item(27,4).Value = 86
# New row 28
$ws.Cells.Item(28,1).Value = "09:28:08"
$ws.Cells.Item(28,2).Value = "11:14"
$ws.Cells.Item(28,3).Value = "215C_LA PLATA"
$ws.Cells.Item(28,4).Value = 106
$ws.Cells.Item(28,5).Value = "L6203"
